$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now holds data previously in row 21
$ws.Range("D2").Value = 44295
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 1800
$ws.Range("M2").Value = 1650
$ws.Range("P2").Value = 550

# Row 3 now holds data previously in row 14
$ws.Range("D3").Value = 44278
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 140
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = 2250
$ws.Range("P3").Value = 750

# Row 4 now holds data previously in row 15
$ws.Range("D4").Value = 44278
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 1800
$ws.Range("M4").Value = 1650
$ws.Range("P4").Value = 550

# Row 5 now holds data previously in row 27
$ws.Range("D5").Value = 44432
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 1200
$ws.Range("L5").Value = 1300
$ws.Range("M5").Value = 1250
$ws.Range("P5").Value = 417

# Row 6 now holds data previously in row 28
$ws.Range("D6").Value = 44432
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 950
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = 975
$ws.Range("P6").Value = 325

# Row 7 now holds data previously in row 13
$ws.Range("D7").Value = 44174
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 500
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = 550
$ws.Range("P7").Value = 183

# Row 8 now holds data previously in row 34
$ws.Range("D8").Value = 44428
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 1500
$ws.Range("L8").Value = 1800
$ws.Range("M8").Value = 1650
$ws.Range("P8").Value = 550

# Row 9 now holds data previously in row 24
$ws.Range("D9").Value = 44398
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 1700
$ws.Range("L9").Value = 1800
$ws.Range("M9").Value = 1750
$ws.Range("P9").Value = 583

# Row 10 now holds data previously in row 7
$ws.Range("D10").Value = 44364
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 1700
$ws.Range("L10").Value = 1800
$ws.Range("M10").Value = 1750
$ws.Range("P10").Value = 583

# Row 11 now holds data previously in row 8
$ws.Range("D11").Value = 44364
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 1400
$ws.Range("L11").Value = 1500
$ws.Range("M11").Value = 1450
$ws.Range("P11").Value = 483

# Row 13 now holds data previously in row 30
$ws.Range("D13").Value = 44302
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 1400
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = 1450
$ws.Range("P13").Value = 483

# Row 14 now holds data previously in row 35
$ws.Range("D14").Value = 44435
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 450
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 1300
$ws.Range("M14").Value = 1194
$ws.Range("P14").Value = 398

# Row 15 now holds data previously in row 36
$ws.Range("D15").Value = 44435
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 950
$ws.Range("L15").Value = 1000
$ws.Range("M15").Value = 975
$ws.Range("P15").Value = 325

# Row 16 now holds data previously in row 19
$ws.Range("D16").Value = 44224
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 1400
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = 1450
$ws.Range("P16").Value = 483

# Row 17 now holds data previously in row 20
$ws.Range("D17").Value = 44224
$ws.Range("I17").Value = "Segunda"
$ws.Range("J17").Value = 160
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 1200
$ws.Range("M17").Value = 1100
$ws.Range("P17").Value = 367

# Row 18 now holds data previously in row 10
$ws.Range("D18").Value = 44333
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 1500
$ws.Range("L18").Value = 1700
$ws.Range("M18").Value = 1600
$ws.Range("P18").Value = 533

# Row 19 now holds data previously in row 16
$ws.Range("D19").Value = 44385
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = 2000
$ws.Range("L19").Value = 2300
$ws.Range("M19").Value = 2150
$ws.Range("P19").Value = 717

# Row 20 now holds data previously in row 31
$ws.Range("D20").Value = 44391
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 250
$ws.Range("K20").Value = 1800
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = 1900
$ws.Range("P20").Value = 633

# Row 21 now holds data previously in row 25
$ws.Range("D21").Value = 44417
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 250
$ws.Range("K21").Value = 1800
$ws.Range("L21").Value = 2000
$ws.Range("M21").Value = 1900
$ws.Range("P21").Value = 633

# Row 22 now holds data previously in row 26
$ws.Range("D22").Value = 44417
$ws.Range("I22").Value = "Segunda"
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 1600
$ws.Range("M22").Value = 1550
$ws.Range("P22").Value = 517

# Row 23 now holds data previously in row 29
$ws.Range("D23").Value = 44342
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 2000
$ws.Range("L23").Value = 2200
$ws.Range("M23").Value = 2100
$ws.Range("P23").Value = 700

# Row 24 now holds data previously in row 2
$ws.Range("D24").Value = 44327
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 200
$ws.Range("K24").Value = 1400
$ws.Range("L24").Value = 1500
$ws.Range("M24").Value = 1450
$ws.Range("P24").Value = 483

# Row 25 now holds data previously in row 18
$ws.Range("D25").Value = 44300
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 250
$ws.Range("K25").Value = 1600
$ws.Range("L25").Value = 1800
$ws.Range("M25").Value = 1700
$ws.Range("P25").Value = 567

# Row 26 now holds data previously in row 32
$ws.Range("D26").Value = 44161
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 200
$ws.Range("K26").Value = 600
$ws.Range("L26").Value = 700
$ws.Range("M26").Value = 650
$ws.Range("P26").Value = 217

# Row 27 now holds data previously in row 33
$ws.Range("D27").Value = 44161
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 250
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 600
$ws.Range("M27").Value = 550
$ws.Range("P27").Value = 183

# Row 28 now holds data previously in row 22
$ws.Range("D28").Value = 44249
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 1500
$ws.Range("L28").Value = 1600
$ws.Range("M28").Value = 1550
$ws.Range("P28").Value = 517

# Row 29 now holds data previously in row 3
$ws.Range("D29").Value = 44280
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 1800
$ws.Range("L29").Value = 2000
$ws.Range("M29").Value = 1900
$ws.Range("P29").Value = 633

# Row 30 now holds data previously in row 4
$ws.Range("D30").Value = 44280
$ws.Range("I30").Value = "Segunda"
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 1400
$ws.Range("L30").Value = 1500
$ws.Range("M30").Value = 1450
$ws.Range("P30").Value = 483

# Row 31 now holds data previously in row 17
$ws.Range("D31").Value = 44447
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 300
$ws.Range("K31").Value = 1100
$ws.Range("L31").Value = 1200
$ws.Range("M31").Value = 1150
$ws.Range("P31").Value = 383

# Row 32 now holds data previously in row 5
$ws.Range("D32").Value = 44270
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 1800
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = 1900
$ws.Range("P32").Value = 633

# Row 33 now holds data previously in row 6
$ws.Range("D33").Value = 44270
$ws.Range("I33").Value = "Segunda"
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 1200
$ws.Range("L33").Value = 1500
$ws.Range("M33").Value = 1350
$ws.Range("P33").Value = 450

# Row 34 now holds data previously in row 11
$ws.Range("D34").Value = 44166
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 250
$ws.Range("K34").Value = 900
$ws.Range("L34").Value = 1000
$ws.Range("M34").Value = 950
$ws.Range("P34").Value = 317

# Row 35 now holds data previously in row 23
$ws.Range("D35").Value = 44397
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 200
$ws.Range("K35").Value = 1400
$ws.Range("L35").Value = 1500
$ws.Range("M35").Value = 1450
$ws.Range("P35").Value = 483

# Row 36 now holds data previously in row 9
$ws.Range("D36").Value = 44306
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 200
$ws.Range("K36").Value = 2400
$ws.Range("L36").Value = 2500
$ws.Range("M36").Value = 2450
$ws.Range("P36").Value = 817
